# search photo command can now search photos for all cars of the set
# Adds a "Bot Comparison" table (G:J, rows 1/3-13) to the log/availabilities
# sheet, comparing this project's bots against TrackPulse VIC / VPT Bot /
# trainbot feature-by-feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New comparison table -------------------------------------------------
# Values are written in the same order the author would have typed them so
# new shared-string entries land at the indices the target file expects.

$ws.Range("H3").Value = "TrackPulse VIC"
$ws.Range("H1").Value = "Bot Comparison"
$ws.Range("I3").Value = "VPT Bot"
$ws.Range("J3").Value = "trainbot"

$ws.Range("G4").Value = "Train Line Status"
$ws.Range("H4").Value = "✅"
$ws.Range("I4").Value = "✅"
$ws.Range("J4").Value = "✅"

$ws.Range("G5").Value = "Bus and Tram Line Status"
$ws.Range("H5").Value = "✅"
$ws.Range("I5").Value = "❎"
$ws.Range("J5").Value = "✅"

$ws.Range("G6").Value = "Station Departures"
$ws.Range("H6").Value = "✅"
$ws.Range("I6").Value = "✅"
$ws.Range("J6").Value = "✅"

$ws.Range("G7").Value = "Train search"
$ws.Range("H7").Value = "✅"
$ws.Range("I7").Value = "❎"
$ws.Range("J7").Value = "✅"

$ws.Range("G8").Value = "Train Photo"
$ws.Range("H8").Value = "✅"
$ws.Range("I8").Value = "❎"
$ws.Range("J8").Value = "❎"

$ws.Range("G9").Value = "Train Runs and Location"
$ws.Range("H9").Value = "✅"
$ws.Range("I9").Value = "❎"
$cell9 = $ws.Range("J9")
$cell9.Value = "⚠ runs only"
$cell9.Characters(2, 10).Font.Size = 8

$ws.Range("G11").Value = "Log feature"
$ws.Range("H11").Value = "✅"
$ws.Range("I11").Value = "❎"
$cell11 = $ws.Range("J11")
$cell11.Value = "⚠ metro trains only"
$cell11.Characters(3, 17).Font.Size = 8

$ws.Range("G12").Value = "Detailed log stats"
$ws.Range("H12").Value = "✅"
$ws.Range("I12").Value = "❎"
$ws.Range("J12").Value = "❎"

$ws.Range("G13").Value = "Games"
$ws.Range("H13").Value = "✅"
$i13 = $ws.Range("I13")
$i13.Value = "❎"
$i13.Font.Color = 0
$ws.Range("J13").Value = "✅"

# Inserted after the fact between "Train Runs and Location" (row 9) and
# "Log feature" (row 11) - hence its shared string lands last.
$ws.Range("G10").Value = "Tram search"
$ws.Range("H10").Value = "✅"
$i10 = $ws.Range("I10")
$i10.Value = "❎"
$i10.Font.Color = 0
$j10 = $ws.Range("J10")
$j10.Value = "❎"
$j10.Font.Color = 0

# --- Cosmetic / view changes ----------------------------------------------

$ws.Columns("A").ColumnWidth = 8.3
$ws.Columns("G").ColumnWidth = 22.5
$ws.Columns("H").ColumnWidth = 13.65

$ws.PageSetup.Orientation = 1

$window = $excel.ActiveWindow
$window.Zoom = 175
$null = $ws.Range("L7").Select()

